$wb = $excel.ActiveWorkbook

# Update column-type values on the "table content col attribute" sheet:
# char[] / char[] / child table -> string / string / string
$ws2 = $wb.Worksheets.Item("table content col attribute")
$ws2.Range("B3").Value = "string"
$ws2.Range("B4").Value = "string"
$ws2.Range("B5").Value = "string"

# Make this sheet the active tab, with B11 selected (matches the diff's
# new tabSelected/selection + workbook activeTab).
$ws2.Activate()
$ws2.Range("B11").Select()
